# Update the "2019-03-05 study list" worksheet:
#  - row 11 (major_depressive_disorder_27479909): PMID column (C) becomes a
#    real number instead of text, and the heritability columns (I-L) are
#    filled in.
#  - row 22 (depression_27089181): heritability columns (I-L) are filled in.
#  - refresh the remembered selection to J16 (cosmetic, matches the saved
#    workbook state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 - major_depressive_disorder_27479909
$ws.Range("C11").Value = 27479909
$ws.Range("I11").Value = 0.05
$ws.Range("J11").Value = 0.3
$ws.Range("K11").Value = 27479909
$ws.Range("L11").Value = 27089181

# Row 22 - depression_27089181
$ws.Range("I22").Value = 0.1
$ws.Range("J22").Value = 0.3
$ws.Range("K22").Value = "no source"
$ws.Range("L22").Value = 27089181

# Saved selection moves to J16
$ws.Range("J16").Select()
